$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows (2-9) with corrected forecast-error values
$ws.Range("B2").Value = -0.04550424379278954
$ws.Range("C2").Value = 0.7456538105835401
$ws.Range("D2").Value = 0.9754175274495978
$ws.Range("E2").Value = 0.9876322835193257
$ws.Range("F2").Value = 1.0238260247887

$ws.Range("B3").Value = -0.2467860954471002
$ws.Range("C3").Value = 0.6225602973041534
$ws.Range("D3").Value = 0.6404695297949983
$ws.Range("E3").Value = 0.8002934023187985
$ws.Range("F3").Value = 0.7923782760124026

$ws.Range("B4").Value = -0.2262603741813949
$ws.Range("C4").Value = 0.5866867510802006
$ws.Range("D4").Value = 0.6375812921102284
$ws.Range("E4").Value = 0.7984868766048873
$ws.Range("F4").Value = 0.7998096604378206

$ws.Range("B5").Value = -0.1929528870395856
$ws.Range("C5").Value = 0.5902456523507005
$ws.Range("D5").Value = 0.6751142996204421
$ws.Range("E5").Value = 0.8216533938470906
$ws.Range("F5").Value = 0.8376585409962011

$ws.Range("B6").Value = -0.1825922693999573
$ws.Range("C6").Value = 0.6080227476365367
$ws.Range("D6").Value = 0.7020402897296509
$ws.Range("E6").Value = 0.8378784456767288
$ws.Range("F6").Value = 0.8619747050200901

$ws.Range("B7").Value = -0.1975866722726227
$ws.Range("C7").Value = 0.6987575864756779
$ws.Range("D7").Value = 0.7662839300513234
$ws.Range("E7").Value = 0.8753764504779207
$ws.Range("F7").Value = 0.9045158188862696
$ws.Range("G7").Value = 9

$ws.Range("B8").Value = -0.2231692882927663
$ws.Range("C8").Value = 0.705494696849475
$ws.Range("D8").Value = 0.8561635644106061
$ws.Range("E8").Value = 0.925291070102055
$ws.Range("F8").Value = 0.9836822860091603
$ws.Range("G8").Value = 6

$ws.Range("B9").Value = -0.3707202268862896
$ws.Range("C9").Value = 0.3707202268862896
$ws.Range("D9").Value = 0.18691442846061
$ws.Range("E9").Value = 0.4323360133745627
$ws.Range("F9").Value = 0.2724360709542369
$ws.Range("G9").Value = 3

# Add new row 10 for quarter Q8
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Q8"
$ws.Range("B10").Value = -0.3121643713505491
$ws.Range("C10").Value = 0.3121643713505491
$ws.Range("D10").Value = 0.09744659474068354
$ws.Range("E10").Value = 0.3121643713505491
$ws.Range("G10").Value = 1
